$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week date range) ---
$ws.Range("A8").Value = "Volume 30   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/11/2023  Through  9/17/2023"

# --- Crime table updates ---

# Cells that change from a numeric 0/1 count to a text "0" placeholder
# (this report represents a zero count as literal text "0", not numeric 0)
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = "0"
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "0"

# Cells that change from the text "0" placeholder to a real numeric count
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("C22").Value = 1
$ws.Range("F22").NumberFormat = "#,##0"
$ws.Range("F22").Value = 1

# Plain numeric value updates (counts and percent changes)
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 3
$ws.Range("I15").Value = 14
$ws.Range("K15").Value = 16.666666666666
$ws.Range("L15").Value = 7.692307692307
$ws.Range("M15").Value = 27.272727272727
$ws.Range("N15").Value = -22.222222222222
$ws.Range("E16").Value = -100
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -14.285714285714
$ws.Range("J16").Value = 82
$ws.Range("K16").Value = 3.658536585365
$ws.Range("M16").Value = -44.078947368421
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 500
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 84.615384615384
$ws.Range("I17").Value = 176
$ws.Range("J17").Value = 159
$ws.Range("K17").Value = 10.691823899371
$ws.Range("L17").Value = 26.618705035971
$ws.Range("M17").Value = 46.666666666666
$ws.Range("D18").Value = 2
$ws.Range("J18").Value = 69
$ws.Range("K18").Value = -44.927536231884
$ws.Range("L18").Value = -25.490196078431
$ws.Range("M18").Value = -79.347826086956
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 150
$ws.Range("F19").Value = 17
$ws.Range("G19").Value = 13
$ws.Range("H19").Value = 30.769230769230
$ws.Range("I19").Value = 150
$ws.Range("J19").Value = 178
$ws.Range("K19").Value = -15.730337078651
$ws.Range("L19").Value = 44.230769230769
$ws.Range("M19").Value = -7.407407407407
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 110
$ws.Range("J20").Value = 112
$ws.Range("K20").Value = -1.785714285714
$ws.Range("L20").Value = 74.603174603174
$ws.Range("M20").Value = 35.802469135802
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 30
$ws.Range("G21").Value = 56
$ws.Range("H21").Value = 10.714285714285
$ws.Range("I21").Value = 577
$ws.Range("J21").Value = 616
$ws.Range("K21").Value = -6.331168831168
$ws.Range("L21").Value = 41.421568627451
$ws.Range("M21").Value = -19.749652294854
$ws.Range("N21").Value = -80.205831903945
$ws.Range("I22").Value = 8
$ws.Range("K22").Value = 33.333333333333
$ws.Range("L22").Value = 33.333333333333
$ws.Range("M22").Value = -20
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 8
$ws.Range("H23").Value = 60
$ws.Range("J23").Value = 75
$ws.Range("K23").Value = 18.666666666666
$ws.Range("L23").Value = 36.923076923076
$ws.Range("M23").Value = 117.073170731707
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = -11.111111111111
$ws.Range("F24").Value = 36
$ws.Range("G24").Value = 39
$ws.Range("H24").Value = -7.692307692307
$ws.Range("I24").Value = 414
$ws.Range("J24").Value = 398
$ws.Range("K24").Value = 4.020100502512
$ws.Range("L24").Value = 43.252595155709
$ws.Range("M24").Value = 3.5
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 16.666666666666
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = -32.352941176470
$ws.Range("I25").Value = 265
$ws.Range("J25").Value = 284
$ws.Range("K25").Value = -6.690140845070
$ws.Range("L25").Value = 22.119815668202
$ws.Range("M25").Value = -22.740524781341
$ws.Range("C26").Value = 1
$ws.Range("F26").Value = 3
$ws.Range("I26").Value = 17
$ws.Range("K26").Value = -15
$ws.Range("L26").Value = -10.526315789473
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 400
$ws.Range("I27").Value = 21
$ws.Range("K27").Value = 5
$ws.Range("L27").Value = 5
$ws.Range("H28").Value = -100
$ws.Range("N28").Value = -50
$ws.Range("H29").Value = -100
$ws.Range("N29").Value = -33.333333333333
